$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 42608.901689814818
$ws.Range("B4").Value = -22
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 53
$ws.Range("E4").Value = 29
$ws.Range("F4").Value = 70
$ws.Range("G4").Value = 18795
$ws.Range("H4").Value = 20605
$ws.Range("I4").Value = 2957
$ws.Range("J4").Value = 369
$ws.Range("K4").Value = 485
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 12
